$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.AddTextbox(1, 537.9414960629921, 67.63637795275591, 3.5999212598425196, 421.67811023622045)
$shp.Name = "TextBox 3"
$shp.Fill.Visible = $false
$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = 1
$shp.TextFrame.TextRange.Text = "dhasudyiusadyiewdsa"
